$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: no faktur (A), nama user (B), total (C), dpp (D), ppn (E)
$data = @(
    @(100082444248363, "RSU PKU MUHAMMADIYAH DELANGGU", 2020000, 1819820, 200180),
    @(100082444248365, "FOCUS INDEPENDEN SCHOOL", 1400000, 1261261, 138739),
    @(100082444248367, "PT Vinsa Mandiri Utama III", 1200000, 1081081, 118919),
    @(100082444248368, "YAYASAN LEMBAGA ELTI GRAMEDIA", 1200000, 1081081, 118919),
    @(100082444248369, "SMA AL AZHAR SYIFA BUDI SURAKARTA", 1200000, 1081081, 118919)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $row++
}

# Column A width, matching Excel's "best fit" autofit width of 12 characters
$ws.Columns.Item(1).ColumnWidth = 11.17

# Set the active selection to B6, matching the diff
$ws.Range("B6").Select() | Out-Null
